$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.058.98"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "3.528.99"
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.26"
$ws.Range("E5").Value = "  +2.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.96"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("D7").Value = "3.526.61"
$ws.Range("E7").Value = "  +3.64%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("E10").Value = "  +3.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.90"
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("E12").Value = "  +3.47%  "
$ws.Range("D13").Value = "4.131.15"
$ws.Range("E13").Value = "  +3.86%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.16"
$ws.Range("E15").Value = "  +4.55%  "
$ws.Range("D16").Value = "3.528.98"
$ws.Range("E16").Value = "  +3.70%  "
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "65.028.71"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.02"
$ws.Range("E19").Value = "  +5.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.83"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.19"
$ws.Range("E21").Value = "  +6.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.36"
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.573"
$ws.Range("E23").Value = "  +5.38%  "
$ws.Range("D24").Value = "3.668.23"
$ws.Range("E24").Value = "  +3.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.87"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000113"
$ws.Range("E27").Value = "  +9.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.68"
$ws.Range("E28").Value = "  +9.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.25"
$ws.Range("E30").Value = "  +4.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.19"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("D32").Value = "3.544.82"
$ws.Range("E32").Value = "  +3.84%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.72"
$ws.Range("E34").Value = "  +4.05%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.36"
$ws.Range("E35").Value = "  +16.87%  "
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "170.07"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("E38").Value = "  +8.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.82"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.95"
$ws.Range("E40").Value = "  +8.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0797"
$ws.Range("E41").Value = "  +7.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.819"
$ws.Range("E42").Value = "  +1.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.54"
$ws.Range("E43").Value = "  +18.54%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.46"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.42"
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.20"
$ws.Range("E47").Value = "  +7.97%  "
$ws.Range("E48").Value = "  +5.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.84"
$ws.Range("E49").Value = "  +7.37%  "
$ws.Range("D50").Value = "2.407.95"
$ws.Range("E50").Value = "  +12.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "301.84"
$ws.Range("E51").Value = "  +11.15%  "
